# GSC export update: append 4 new daily rows (2025-11-07 .. 2025-11-10)
# to the "Chart" worksheet, matching the pattern of the preceding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$dates = @("2025-11-07", "2025-11-08", "2025-11-09", "2025-11-10")
$noVideoIndexed = 24.0
$videoIndexed = 0.0
$impressions = 0.0

$startRow = 35
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i

    # Write the date as a genuine text value (matching the existing date
    # cells, which are shared-string text rather than date serials).
    # A direct .Value assignment of a date-shaped string gets auto-converted
    # to a date serial number by Excel, so instead we enter it as a literal
    # text formula and then collapse the formula down to its cached value
    # via copy / paste-values. That keeps the cell's style untouched (same
    # default style as every other cell) instead of stamping a new number
    # format onto it.
    $cellA = $ws.Cells.Item($r, 1)
    $cellA.Formula = '="' + $dates[$i] + '"'
    $cellA.Copy()
    $cellA.PasteSpecial(-4163)

    $ws.Cells.Item($r, 2).Value = $noVideoIndexed
    $ws.Cells.Item($r, 3).Value = $videoIndexed
    $ws.Cells.Item($r, 4).Value = $impressions
}

$excel.CutCopyMode = 0
